$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44699
$ws.Cells.Item(2, 11).Value = 'Mankaki'
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 250
$ws.Cells.Item(2, 14).Value = 29000
$ws.Cells.Item(2, 15).Value = 30000
$ws.Cells.Item(2, 16).Value = 29500
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 1639

# Row 3
$ws.Cells.Item(3, 4).Value = 45071
$ws.Cells.Item(3, 11).Value = 'Fuyu'
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 110
$ws.Cells.Item(3, 14).Value = 23000
$ws.Cells.Item(3, 15).Value = 24000
$ws.Cells.Item(3, 16).Value = 23455
$ws.Cells.Item(3, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 19).Value = 1303

# Row 4
$ws.Cells.Item(4, 4).Value = 44355
$ws.Cells.Item(4, 11).Value = 'Mankaki'
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 13).Value = 270
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 20500
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 1139

# Row 5
$ws.Cells.Item(5, 4).Value = 44301
$ws.Cells.Item(5, 11).Value = 'Hachiya'
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 250
$ws.Cells.Item(5, 14).Value = 20000
$ws.Cells.Item(5, 15).Value = 21000
$ws.Cells.Item(5, 16).Value = 20500
$ws.Cells.Item(5, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 19).Value = 1139

# Row 6
$ws.Cells.Item(6, 4).Value = 44305
$ws.Cells.Item(6, 11).Value = 'Mankaki'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 250
$ws.Cells.Item(6, 14).Value = 24000
$ws.Cells.Item(6, 15).Value = 25000
$ws.Cells.Item(6, 16).Value = 24500
$ws.Cells.Item(6, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 19).Value = 1361

# Row 7
$ws.Cells.Item(7, 4).Value = 44313
$ws.Cells.Item(7, 11).Value = 'Mankaki'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 270
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7, 19).Value = 1194

# Row 8
$ws.Cells.Item(8, 4).Value = 45043
$ws.Cells.Item(8, 11).Value = 'Fuyu'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 300
$ws.Cells.Item(8, 14).Value = 25000
$ws.Cells.Item(8, 15).Value = 26000
$ws.Cells.Item(8, 16).Value = 25500
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1417
